# Auto-generated Excel COM-interop script
# Applies numeric updates to the Leve profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as described by the upstream diff (scheduled market-price data refresh).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 46: Always Have an Exit Plan | Poisoning Potion
$ws.Cells.Item(46, 8).Value = 5607.3335  # H46: 20004080 -> 5607.3335
$ws.Cells.Item(46, 9).Value = 7400  # I46: 5200 -> 7400
$ws.Cells.Item(46, 10).Value = 4711  # J46: 33336666 -> 4711
$ws.Cells.Item(46, 11).Value = 22200  # K46: 15600 -> 22200
$ws.Cells.Item(46, 12).Value = 14133  # L46: 100009998 -> 14133
$ws.Cells.Item(46, 13).Value = -22081  # M46: -15481 -> -22081
$ws.Cells.Item(46, 14).Value = -14371  # N46: -100010236 -> -14371

# Row 60: Make Up Your Mind or Else | Potent Poisoning Potion
$ws.Cells.Item(60, 8).Value = 5607.3335  # H60: 20004080 -> 5607.3335
$ws.Cells.Item(60, 9).Value = 7400  # I60: 5200 -> 7400
$ws.Cells.Item(60, 10).Value = 4711  # J60: 33336666 -> 4711
$ws.Cells.Item(60, 11).Value = 22200  # K60: 15600 -> 22200
$ws.Cells.Item(60, 12).Value = 14133  # L60: 100009998 -> 14133
$ws.Cells.Item(60, 13).Value = -21716  # M60: -15116 -> -21716
$ws.Cells.Item(60, 14).Value = -15101  # N60: -100010966 -> -15101

# Row 98: The Dotted Line | Enchanted Durium Ink
$ws.Cells.Item(98, 8).Value = 1408.375  # H98: 1110.9048 -> 1408.375
$ws.Cells.Item(98, 9).Value = 1354.5  # I98: 1042.0714 -> 1354.5
$ws.Cells.Item(98, 10).Value = 1570  # J98: 1248.5714 -> 1570
$ws.Cells.Item(98, 11).Value = 1354.5  # K98: 1042.0714 -> 1354.5
$ws.Cells.Item(98, 12).Value = 1570  # L98: 1248.5714 -> 1570
$ws.Cells.Item(98, 13).Value = 143.5  # M98: 455.9286 -> 143.5
$ws.Cells.Item(98, 14).Value = -4566  # N98: -4244.5714 -> -4566

# Row 113: Amaro Kart | Starch Glue
$ws.Cells.Item(113, 8).Value = 3128.9167  # H113: 3187.3333 -> 3128.9167
$ws.Cells.Item(113, 9).Value = 2176.25  # I113: 2233.3333 -> 2176.25
$ws.Cells.Item(113, 10).Value = 3605.25  # J113: 3505.3333 -> 3605.25
$ws.Cells.Item(113, 11).Value = 2176.25  # K113: 2233.3333 -> 2176.25
$ws.Cells.Item(113, 12).Value = 3605.25  # L113: 3505.3333 -> 3605.25
$ws.Cells.Item(113, 13).Value = 1077.75  # M113: 1020.6667 -> 1077.75
$ws.Cells.Item(113, 14).Value = -10113.25  # N113: -10013.3333 -> -10113.25

# Row 121: Mindful Medicine | Tincture of Mind
$ws.Cells.Item(121, 8).Value = 3900  # H121: 3025 -> 3900
$ws.Cells.Item(121, 9).Value = 1800  # I121: 1333.3334 -> 1800
$ws.Cells.Item(121, 11).Value = 5400  # K121: 4000.0002 -> 5400
$ws.Cells.Item(121, 13).Value = -3653  # M121: -2253.0002 -> -3653

# Row 122: Wishful Inking | Enchanted High Durium Ink
$ws.Cells.Item(122, 8).Value = 1408.375  # H122: 1110.9048 -> 1408.375
$ws.Cells.Item(122, 9).Value = 1354.5  # I122: 1042.0714 -> 1354.5
$ws.Cells.Item(122, 10).Value = 1570  # J122: 1248.5714 -> 1570
$ws.Cells.Item(122, 11).Value = 4063.5  # K122: 3126.2142 -> 4063.5
$ws.Cells.Item(122, 12).Value = 4710  # L122: 3745.7142 -> 4710
$ws.Cells.Item(122, 13).Value = -1613.5  # M122: -676.2142000000003 -> -1613.5
$ws.Cells.Item(122, 14).Value = -9610  # N122: -8645.7142 -> -9610

# Row 125: Body over Mind | Grade 5 Dexterity Alkahest
$ws.Cells.Item(125, 8).Value = 11645.333  # H125: 11701.777 -> 11645.333
$ws.Cells.Item(125, 9).Value = 476  # I125: 500 -> 476
$ws.Cells.Item(125, 10).Value = 14836.571  # J125: 13102 -> 14836.571
$ws.Cells.Item(125, 11).Value = 4284  # K125: 4500 -> 4284
$ws.Cells.Item(125, 12).Value = 133529.139  # L125: 117918 -> 133529.139
$ws.Cells.Item(125, 13).Value = -1824  # M125: -2040 -> -1824
$ws.Cells.Item(125, 14).Value = -138449.139  # N125: -122838 -> -138449.139

# Row 131: Mindful Study | Grade 5 Tincture of Mind
$ws.Cells.Item(131, 8).Value = 4951.4707  # H131: 3262.5186 -> 4951.4707
$ws.Cells.Item(131, 9).Value = 822.5  # I131: 200.25 -> 822.5
$ws.Cells.Item(131, 10).Value = 6221.923  # J131: 5712.3335 -> 6221.923
$ws.Cells.Item(131, 11).Value = 2467.5  # K131: 600.75 -> 2467.5
$ws.Cells.Item(131, 12).Value = 18665.769  # L131: 17137.0005 -> 18665.769
$ws.Cells.Item(131, 13).Value = 2572.5  # M131: 4439.25 -> 2572.5
$ws.Cells.Item(131, 14).Value = -28745.769  # N131: -27217.0005 -> -28745.769

# Row 141: Remedy for Reason | Grade 1 Gemdraught of Mind
$ws.Cells.Item(141, 8).Value = 4331.316  # H141: 4789.4116 -> 4331.316
$ws.Cells.Item(141, 9).Value = 1798.8889  # I141: 2664.375 -> 1798.8889
$ws.Cells.Item(141, 10).Value = 6610.5  # J141: 6678.3335 -> 6610.5
$ws.Cells.Item(141, 11).Value = 5396.6667  # K141: 7993.125 -> 5396.6667
$ws.Cells.Item(141, 12).Value = 19831.5  # L141: 20035.0005 -> 19831.5
$ws.Cells.Item(141, 13).Value = -216.6666999999998  # M141: -2813.125 -> -216.6666999999998
$ws.Cells.Item(141, 14).Value = -30191.5  # N141: -30395.0005 -> -30191.5

$ws = $wb.Worksheets.Item("ARM")
# Row 8: You've Got Mail | Bronze Haubergeon
$ws.Cells.Item(8, 8).Value = 502500  # H8: 1000000 -> 502500
$ws.Cells.Item(8, 10).Value = 5000  # J8: 0 -> 5000
$ws.Cells.Item(8, 12).Value = 5000  # L8: 0 -> 5000
$ws.Cells.Item(8, 14).Value = -5288  # N8: None -> -5288

# Row 92: Mail It In | High Steel Scale Mail of Fending
$ws.Cells.Item(92, 8).Value = 33000  # H92: 25000 -> 33000
$ws.Cells.Item(92, 10).Value = 33000  # J92: 25000 -> 33000
$ws.Cells.Item(92, 12).Value = 33000  # L92: 25000 -> 33000
$ws.Cells.Item(92, 14).Value = -37992  # N92: -29992 -> -37992

# Row 96: The Gauntlet Is Cast | High Steel Gauntlets of Fending
$ws.Cells.Item(96, 8).Value = 0  # H96: 17312 -> 0
$ws.Cells.Item(96, 9).Value = 0  # I96: 17312 -> 0
$ws.Cells.Item(96, 11).Value = 0  # K96: 17312 -> 0
$ws.Cells.Item(96, 13).ClearContents()  # M96: remove (was -14566)

# Row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
$ws.Cells.Item(102, 8).Value = 3089.4736  # H102: 3255.8823 -> 3089.4736
$ws.Cells.Item(102, 10).Value = 3310  # J102: 4400 -> 3310
$ws.Cells.Item(102, 12).Value = 3310  # L102: 4400 -> 3310
$ws.Cells.Item(102, 14).Value = -6554  # N102: -7644 -> -6554

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Cells.Item(132, 8).Value = 2554.5518  # H132: 2739.5557 -> 2554.5518
$ws.Cells.Item(132, 9).Value = 1753.5  # I132: 1997.4117 -> 1753.5
$ws.Cells.Item(132, 10).Value = 4334.6665  # J132: 4001.2 -> 4334.6665
$ws.Cells.Item(132, 11).Value = 5260.5  # K132: 5992.2351 -> 5260.5
$ws.Cells.Item(132, 12).Value = 13003.9995  # L132: 12003.6 -> 13003.9995
$ws.Cells.Item(132, 13).Value = -2730.5  # M132: -3462.2351 -> -2730.5
$ws.Cells.Item(132, 14).Value = -18063.9995  # N132: -17063.6 -> -18063.9995

$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run | Iron Rivets
$ws.Cells.Item(22, 8).Value = 0  # H22: 400 -> 0
$ws.Cells.Item(22, 10).Value = 0  # J22: 400 -> 0
$ws.Cells.Item(22, 12).Value = 0  # L22: 400 -> 0
$ws.Cells.Item(22, 14).ClearContents()  # N22: remove (was -746)

# Row 37: That's Some Fine Grinding | Initiate's Mortar
$ws.Cells.Item(37, 8).Value = 3902  # H37: 7533.3335 -> 3902
$ws.Cells.Item(37, 9).Value = 1863  # I37: 3000 -> 1863
$ws.Cells.Item(37, 10).Value = 7980  # J37: 9800 -> 7980
$ws.Cells.Item(37, 11).Value = 1863  # K37: 3000 -> 1863
$ws.Cells.Item(37, 12).Value = 7980  # L37: 9800 -> 7980
$ws.Cells.Item(37, 13).Value = -1726  # M37: -2863 -> -1726
$ws.Cells.Item(37, 14).Value = -8254  # N37: -10074 -> -8254

$ws = $wb.Worksheets.Item("CRP")
# Row 17: Say It with Spears | Feathered Harpoon
$ws.Cells.Item(17, 8).Value = 20400  # H17: 23850 -> 20400
$ws.Cells.Item(17, 10).Value = 25800  # J17: 26800 -> 25800
$ws.Cells.Item(17, 12).Value = 25800  # L17: 26800 -> 25800
$ws.Cells.Item(17, 14).Value = -26148  # N17: -27148 -> -26148

# Row 54: The Turning Point | Garnet Grinding Wheel
$ws.Cells.Item(54, 8).Value = 30297.334  # H54: 40046 -> 30297.334
$ws.Cells.Item(54, 10).Value = 30297.334  # J54: 40046 -> 30297.334
$ws.Cells.Item(54, 12).Value = 30297.334  # L54: 40046 -> 30297.334
$ws.Cells.Item(54, 14).Value = -31613.334  # N54: -41362 -> -31613.334

# Row 55: Ready for a Rematch | Mythril Lance
$ws.Cells.Item(55, 8).Value = 12416.5  # H55: 11000 -> 12416.5
$ws.Cells.Item(55, 9).Value = 13000  # I55: 7000 -> 13000
$ws.Cells.Item(55, 10).Value = 12222  # J55: 15000 -> 12222
$ws.Cells.Item(55, 11).Value = 13000  # K55: 7000 -> 13000
$ws.Cells.Item(55, 12).Value = 12222  # L55: 15000 -> 12222
$ws.Cells.Item(55, 13).Value = -12685  # M55: -6685 -> -12685
$ws.Cells.Item(55, 14).Value = -12852  # N55: -15630 -> -12852

# Row 59: Bow Down to Magic | Crab Bow
$ws.Cells.Item(59, 8).Value = 20121  # H59: 19642.9 -> 20121
$ws.Cells.Item(59, 10).Value = 20121  # J59: 19642.9 -> 20121
$ws.Cells.Item(59, 12).Value = 20121  # L59: 19642.9 -> 20121
$ws.Cells.Item(59, 14).Value = -22411  # N59: -21932.9 -> -22411

# Row 60: Bowing to Greater Power | Yew Longbow
$ws.Cells.Item(60, 8).Value = 10302.115  # H60: 14289.818 -> 10302.115
$ws.Cells.Item(60, 9).Value = 2050  # I60: 0 -> 2050
$ws.Cells.Item(60, 10).Value = 10989.792  # J60: 14289.818 -> 10989.792
$ws.Cells.Item(60, 11).Value = 2050  # K60: 0 -> 2050
$ws.Cells.Item(60, 12).Value = 10989.792  # L60: 14289.818 -> 10989.792
$ws.Cells.Item(60, 13).Value = -1539  # M60: None -> -1539
$ws.Cells.Item(60, 14).Value = -12011.792  # N60: -15311.818 -> -12011.792

# Row 111: Taking Aim | Applewood Longbow
$ws.Cells.Item(111, 8).Value = 79800  # H111: 70000 -> 79800
$ws.Cells.Item(111, 10).Value = 79800  # J111: 70000 -> 79800
$ws.Cells.Item(111, 12).Value = 79800  # L111: 70000 -> 79800
$ws.Cells.Item(111, 14).Value = -87980  # N111: -78180 -> -87980

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Cells.Item(132, 8).Value = 2548.524  # H132: 2747 -> 2548.524
$ws.Cells.Item(132, 9).Value = 2334.1482  # I132: 2548.75 -> 2334.1482
$ws.Cells.Item(132, 10).Value = 2934.4  # J132: 3086.8572 -> 2934.4
$ws.Cells.Item(132, 11).Value = 7002.444600000001  # K132: 7646.25 -> 7002.444600000001
$ws.Cells.Item(132, 12).Value = 8803.200000000001  # L132: 9260.571599999999 -> 8803.200000000001
$ws.Cells.Item(132, 13).Value = -4472.444600000001  # M132: -5116.25 -> -4472.444600000001
$ws.Cells.Item(132, 14).Value = -13863.2  # N132: -14320.5716 -> -13863.2

$ws = $wb.Worksheets.Item("CUL")
# Row 87: Soup That Eats Like a Knight | Clam Chowder
$ws.Cells.Item(87, 8).Value = 7743  # H87: 6885.4287 -> 7743
$ws.Cells.Item(87, 9).Value = 1000  # I87: 3000 -> 1000
$ws.Cells.Item(87, 10).Value = 8192.532999999999  # J87: 9799.5 -> 8192.532999999999
$ws.Cells.Item(87, 11).Value = 3000  # K87: 9000 -> 3000
$ws.Cells.Item(87, 12).Value = 24577.599  # L87: 29398.5 -> 24577.599
$ws.Cells.Item(87, 13).Value = -1752  # M87: -7752 -> -1752
$ws.Cells.Item(87, 14).Value = -27073.599  # N87: -31894.5 -> -27073.599

# Row 90: Like Ma Used to Make (L) | Clam Chowder
$ws.Cells.Item(90, 8).Value = 7743  # H90: 6885.4287 -> 7743
$ws.Cells.Item(90, 9).Value = 1000  # I90: 3000 -> 1000
$ws.Cells.Item(90, 10).Value = 8192.532999999999  # J90: 9799.5 -> 8192.532999999999
$ws.Cells.Item(90, 11).Value = 9000  # K90: 27000 -> 9000
$ws.Cells.Item(90, 12).Value = 73732.79699999999  # L90: 88195.5 -> 73732.79699999999
$ws.Cells.Item(90, 13).Value = -2760  # M90: -20760 -> -2760
$ws.Cells.Item(90, 14).Value = -86212.79699999999  # N90: -100675.5 -> -86212.79699999999

# Row 113: Can't Eat Just One | Night Vinegar
$ws.Cells.Item(113, 8).Value = 685.3585  # H113: 678.25 -> 685.3585
$ws.Cells.Item(113, 9).Value = 696.6739  # I113: 691.4681 -> 696.6739
$ws.Cells.Item(113, 10).Value = 611  # J113: 609.2222 -> 611
$ws.Cells.Item(113, 11).Value = 2090.0217  # K113: 2074.4043 -> 2090.0217
$ws.Cells.Item(113, 12).Value = 1833  # L113: 1827.6666 -> 1833
$ws.Cells.Item(113, 13).Value = 79.97829999999976  # M113: 95.59569999999985 -> 79.97829999999976
$ws.Cells.Item(113, 14).Value = -6173  # N113: -6167.6666 -> -6173

# Row 132: More Mezcal | Cooking Mezcal
$ws.Cells.Item(132, 8).Value = 2208.6667  # H132: 2580.5 -> 2208.6667
$ws.Cells.Item(132, 9).Value = 3049.75  # I132: 5400 -> 3049.75
$ws.Cells.Item(132, 10).Value = 1788.125  # J132: 1875.625 -> 1788.125
$ws.Cells.Item(132, 11).Value = 27447.75  # K132: 48600 -> 27447.75
$ws.Cells.Item(132, 12).Value = 16093.125  # L132: 16880.625 -> 16093.125
$ws.Cells.Item(132, 13).Value = -24917.75  # M132: -46070 -> -24917.75
$ws.Cells.Item(132, 14).Value = -21153.125  # N132: -21940.625 -> -21153.125

# Row 140: Sweet, Sweet Bean Juice | Mesquite Juice
$ws.Cells.Item(140, 8).Value = 2205.6775  # H140: 2612.2856 -> 2205.6775
$ws.Cells.Item(140, 9).Value = 1601  # I140: 1818.9333 -> 1601
$ws.Cells.Item(140, 10).Value = 3305.0908  # J140: 3527.6924 -> 3305.0908
$ws.Cells.Item(140, 11).Value = 4803  # K140: 5456.7999 -> 4803
$ws.Cells.Item(140, 12).Value = 9915.2724  # L140: 10583.0772 -> 9915.2724
$ws.Cells.Item(140, 13).Value = 377  # M140: -276.7999 -> 377
$ws.Cells.Item(140, 14).Value = -20275.2724  # N140: -20943.0772 -> -20275.2724

$ws = $wb.Worksheets.Item("GSM")
# Row 47: Wear Your Patriotic Pin | Peridot Choker
$ws.Cells.Item(47, 8).Value = 19000  # H47: 20015.5 -> 19000
$ws.Cells.Item(47, 10).Value = 19000  # J47: 20015.5 -> 19000
$ws.Cells.Item(47, 12).Value = 19000  # L47: 20015.5 -> 19000
$ws.Cells.Item(47, 14).Value = -20136  # N47: -21151.5 -> -20136

# Row 122: Awarding Academic Excellence | Ametrine
$ws.Cells.Item(122, 8).Value = 17999.666  # H122: 5319.923 -> 17999.666
$ws.Cells.Item(122, 9).Value = 49999  # I122: 6335.1 -> 49999
$ws.Cells.Item(122, 10).Value = 2000  # J122: 1936 -> 2000
$ws.Cells.Item(122, 11).Value = 149997  # K122: 19005.3 -> 149997
$ws.Cells.Item(122, 12).Value = 6000  # L122: 5808 -> 6000
$ws.Cells.Item(122, 13).Value = -147547  # M122: -16555.3 -> -147547
$ws.Cells.Item(122, 14).Value = -10900  # N122: -10708 -> -10900

$ws = $wb.Worksheets.Item("LTW")
# Row 35: No Risk, No Reward | Toadskin Cesti
$ws.Cells.Item(35, 8).Value = 18739.625  # H35: 23946.154 -> 18739.625
$ws.Cells.Item(35, 9).Value = 7969.8  # I35: 3250 -> 7969.8
$ws.Cells.Item(35, 10).Value = 23635  # J35: 27709.092 -> 23635
$ws.Cells.Item(35, 11).Value = 7969.8  # K35: 3250 -> 7969.8
$ws.Cells.Item(35, 12).Value = 23635  # L35: 27709.092 -> 23635
$ws.Cells.Item(35, 13).Value = -7633.8  # M35: -2914 -> -7633.8
$ws.Cells.Item(35, 14).Value = -24307  # N35: -28381.092 -> -24307

# Row 61: Spelling Me Softly | Raptor Leather
$ws.Cells.Item(61, 8).Value = 465397.03  # H61: 511796.5 -> 465397.03
$ws.Cells.Item(61, 9).Value = 17524.615  # I61: 15268 -> 17524.615
$ws.Cells.Item(61, 10).Value = 1112323.9  # J61: 2001382 -> 1112323.9
$ws.Cells.Item(61, 11).Value = 17524.615  # K61: 15268 -> 17524.615
$ws.Cells.Item(61, 12).Value = 1112323.9  # L61: 2001382 -> 1112323.9
$ws.Cells.Item(61, 13).Value = -17322.615  # M61: -15066 -> -17322.615
$ws.Cells.Item(61, 14).Value = -1112727.9  # N61: -2001786 -> -1112727.9

# Row 113: Peace in Rest | Atrociraptor Leather
$ws.Cells.Item(113, 8).Value = 465397.03  # H113: 511796.5 -> 465397.03
$ws.Cells.Item(113, 9).Value = 17524.615  # I113: 15268 -> 17524.615
$ws.Cells.Item(113, 10).Value = 1112323.9  # J113: 2001382 -> 1112323.9
$ws.Cells.Item(113, 11).Value = 17524.615  # K113: 15268 -> 17524.615
$ws.Cells.Item(113, 12).Value = 1112323.9  # L113: 2001382 -> 1112323.9
$ws.Cells.Item(113, 13).Value = -15354.615  # M113: -13098 -> -15354.615
$ws.Cells.Item(113, 14).Value = -1116663.9  # N113: -2005722 -> -1116663.9

# Row 137: Lending Artisans a Hand | Br'aaxskin Halfgloves of Crafting
$ws.Cells.Item(137, 8).Value = 37623.625  # H137: 49997.5 -> 37623.625
$ws.Cells.Item(137, 10).Value = 37623.625  # J137: 49997.5 -> 37623.625
$ws.Cells.Item(137, 12).Value = 37623.625  # L137: 49997.5 -> 37623.625
$ws.Cells.Item(137, 14).Value = -47823.625  # N137: -60197.5 -> -47823.625

$ws = $wb.Worksheets.Item("WVR")
# Row 45: Private Concerns | Linen Trousers
$ws.Cells.Item(45, 8).Value = 4692  # H45: 4643.75 -> 4692
$ws.Cells.Item(45, 10).Value = 4692  # J45: 4643.75 -> 4692
$ws.Cells.Item(45, 12).Value = 4692  # L45: 4643.75 -> 4692
$ws.Cells.Item(45, 14).Value = -5674  # N45: -5625.75 -> -5674

# Row 47: The Wages of Sin | Linen Coatee of Crafting
$ws.Cells.Item(47, 8).Value = 175000  # H47: 132666.67 -> 175000
$ws.Cells.Item(47, 10).Value = 175000  # J47: 132666.67 -> 175000
$ws.Cells.Item(47, 12).Value = 175000  # L47: 132666.67 -> 175000
$ws.Cells.Item(47, 14).Value = -176144  # N47: -133810.67 -> -176144

# Row 61: Bundle Up, It's Odd out There | Woolen Deerstalker
$ws.Cells.Item(61, 8).Value = 7055.3335  # H61: 8222.5 -> 7055.3335
$ws.Cells.Item(61, 9).Value = 3590.1667  # I61: 4995 -> 3590.1667
$ws.Cells.Item(61, 10).Value = 13985.667  # J61: 11450 -> 13985.667
$ws.Cells.Item(61, 11).Value = 3590.1667  # K61: 4995 -> 3590.1667
$ws.Cells.Item(61, 12).Value = 13985.667  # L61: 11450 -> 13985.667
$ws.Cells.Item(61, 13).Value = -3298.1667  # M61: -4703 -> -3298.1667
$ws.Cells.Item(61, 14).Value = -14569.667  # N61: -12034 -> -14569.667

# Row 96: Skills on Display | Ruby Cotton Cloth
$ws.Cells.Item(96, 8).Value = 1141.75  # H96: 1068.8 -> 1141.75
$ws.Cells.Item(96, 9).Value = 834  # I96: 815 -> 834
$ws.Cells.Item(96, 11).Value = 834  # K96: 815 -> 834
$ws.Cells.Item(96, 13).Value = 539  # M96: 558 -> 539
